$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $val) {
    $fmt = $cell.NumberFormat
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.NumberFormat = $fmt
}

Set-TextValue $ws.Range('D2') '71.921.06'
Set-TextValue $ws.Range('E2') '  -0.56%  '
Set-TextValue $ws.Range('D3') '4.018.74'
Set-TextValue $ws.Range('E3') '  -0.85%  '
Set-TextValue $ws.Range('E4') '  +0.03%  '
Set-TextValue $ws.Range('D5') '534.04'
Set-TextValue $ws.Range('E5') '  +1.36%  '
Set-TextValue $ws.Range('D6') '150.25'
Set-TextValue $ws.Range('E6') '  -0.22%  '
Set-TextValue $ws.Range('D7') '4.014.96'
Set-TextValue $ws.Range('E7') '  -0.65%  '
Set-TextValue $ws.Range('E8') '  -2.73%  '
Set-TextValue $ws.Range('E9') '  -0.04%  '
Set-TextValue $ws.Range('D10') '0.753'
Set-TextValue $ws.Range('E10') '  -2.10%  '
Set-TextValue $ws.Range('D11') '0.171'
Set-TextValue $ws.Range('E11') '  -3.64%  '
Set-TextValue $ws.Range('D12') '54.06'
Set-TextValue $ws.Range('E12') '  +7.47%  '
Set-TextValue $ws.Range('E13') '  -3.17%  '
Set-TextValue $ws.Range('D14') '10.79'
Set-TextValue $ws.Range('E14') '  -2.10%  '
Set-TextValue $ws.Range('D15') '4.658.52'
Set-TextValue $ws.Range('E15') '  -0.92%  '
Set-TextValue $ws.Range('D16') '4.022.78'
Set-TextValue $ws.Range('E16') '  -0.43%  '
Set-TextValue $ws.Range('D17') '14.14'
Set-TextValue $ws.Range('E17') '  -1.76%  '
Set-TextValue $ws.Range('E18') '  -1.58%  '
Set-TextValue $ws.Range('E19') '  -3.61%  '
Set-TextValue $ws.Range('E20') '  -1.57%  '
Set-TextValue $ws.Range('D21') '71.922.61'
Set-TextValue $ws.Range('E21') '  -0.38%  '
Set-TextValue $ws.Range('D22') '432.85'
Set-TextValue $ws.Range('E22') '  -1.28%  '
Set-TextValue $ws.Range('D23') '98.33'
Set-TextValue $ws.Range('E23') '  -2.88%  '
Set-TextValue $ws.Range('D24') '3.60'
Set-TextValue $ws.Range('E24') '  -1.51%  '
Set-TextValue $ws.Range('D25') '14.71'
Set-TextValue $ws.Range('E25') '  -2.29%  '
Set-TextValue $ws.Range('D26') '4.21'
Set-TextValue $ws.Range('E26') '  +0.08%  '
Set-TextValue $ws.Range('D27') '4.36'
Set-TextValue $ws.Range('E27') '  +29.48%  '
Set-TextValue $ws.Range('D28') '11.45'
Set-TextValue $ws.Range('E28') '  -1.50%  '
Set-TextValue $ws.Range('D29') '10.79'
Set-TextValue $ws.Range('E29') '  -2.59%  '
Set-TextValue $ws.Range('D30') '5.94'
Set-TextValue $ws.Range('E30') '  +1.75%  '
Set-TextValue $ws.Range('D31') '36.96'
Set-TextValue $ws.Range('E31') '  -1.71%  '
Set-TextValue $ws.Range('D32') '8.26'
Set-TextValue $ws.Range('E32') '  +21.43%  '
Set-TextValue $ws.Range('E33') '  +1.81%  '
Set-TextValue $ws.Range('D34') '50.37'
Set-TextValue $ws.Range('E34') '  +18.11%  '
Set-TextValue $ws.Range('D35') '13.56'
Set-TextValue $ws.Range('E35') '  -1.21%  '
Set-TextValue $ws.Range('D36') '675.42'
Set-TextValue $ws.Range('E36') '  +0.13%  '
Set-TextValue $ws.Range('D37') '67.63'
Set-TextValue $ws.Range('E37') '  +1.32%  '
Set-TextValue $ws.Range('D38') '0.454'
Set-TextValue $ws.Range('E38') '  +2.95%  '
Set-TextValue $ws.Range('D39') '0.0₃0828'
Set-TextValue $ws.Range('E39') '  -5.03%  '
Set-TextValue $ws.Range('D40') '0.149'
Set-TextValue $ws.Range('E40') '  -5.46%  '
Set-TextValue $ws.Range('D41') '3.42'
Set-TextValue $ws.Range('E41') '  +7.80%  '
Set-TextValue $ws.Range('D42') '3.37'
Set-TextValue $ws.Range('E42') '  -3.19%  '
Set-TextValue $ws.Range('D44') '11.02'
Set-TextValue $ws.Range('E44') '  +15.56%  '
Set-TextValue $ws.Range('D45') '0.0493'
Set-TextValue $ws.Range('E45') '  -3.14%  '
Set-TextValue $ws.Range('D46') '1.00'
Set-TextValue $ws.Range('E46') '  +0.08%  '
Set-TextValue $ws.Range('E47') '  -3.58%  '
Set-TextValue $ws.Range('D48') '2.64'
Set-TextValue $ws.Range('E48') '  -5.55%  '
Set-TextValue $ws.Range('E49') '  -0.22%  '
Set-TextValue $ws.Range('D50') '3.10'
Set-TextValue $ws.Range('E50') '  -0.47%  '
Set-TextValue $ws.Range('D51') '2.848.52'
Set-TextValue $ws.Range('E51') '  +8.75%  '
